$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bulk-write rows 2-52 (all columns) with corrected values
$data = New-Object 'object[,]' 51,5
$data[0,0] = 39400
$data[0,1] = 2007
$data[0,2] = 4.930115226412335
$data[0,3] = 2008
$data[0,4] = 0.3630458632513767
$data[1,0] = 39583
$data[1,1] = 2008
$data[1,2] = 5.021907707863549
$data[1,3] = 2009
$data[1,4] = 19.19812743658083
$data[2,0] = 39765
$data[2,1] = 2008
$data[2,2] = 1.457587285166628
$data[2,3] = 2009
$data[2,4] = 0.507956838644974
$data[3,0] = 39948
$data[3,1] = 2009
$data[3,2] = -5.592633745595466
$data[3,3] = 2010
$data[3,4] = -9.964084247724703
$data[4,0] = 40130
$data[4,1] = 2009
$data[4,2] = -0.9140166223623569
$data[4,3] = 2010
$data[4,4] = 5.6395352704941
$data[5,0] = 40310
$data[5,1] = 2010
$data[5,2] = -5.440152375872254
$data[5,3] = 2011
$data[5,4] = -14.43639438706738
$data[6,0] = 40494
$data[6,1] = 2010
$data[6,2] = 2.585942866987878
$data[6,3] = 2011
$data[6,4] = 1.724360951547554
$data[7,0] = 40676
$data[7,1] = 2011
$data[7,2] = 9.349082908138451
$data[7,3] = 2012
$data[7,4] = 27.15801420548429
$data[8,0] = 40862
$data[8,1] = 2011
$data[8,2] = 4.253963781362402
$data[8,3] = 2012
$data[8,4] = -0.6955733540840336
$data[9,0] = 41044
$data[9,1] = 2012
$data[9,2] = 0.5389546843750148
$data[9,3] = 2013
$data[9,4] = -5.27893918837793
$data[10,0] = 41228
$data[10,1] = 2012
$data[10,2] = 1.752870900283909
$data[10,3] = 2013
$data[10,4] = 4.300339264728548
$data[11,0] = 41409
$data[11,1] = 2013
$data[11,2] = -4.232836797447693
$data[11,3] = 2014
$data[11,4] = -8.0930759205322
$data[12,0] = 41592
$data[12,1] = 2013
$data[12,2] = -1.479696720105139
$data[12,3] = 2014
$data[12,4] = 8.296896928314457
$data[13,0] = 41774
$data[13,1] = 2014
$data[13,2] = 7.942828065321739
$data[13,3] = 2015
$data[13,4] = 15.37760125310905
$data[14,0] = 41957
$data[14,1] = 2014
$data[14,2] = 3.900127535411246
$data[14,3] = 2015
$data[14,4] = -2.092856741436233
$data[15,0] = 42137
$data[15,1] = 2015
$data[15,2] = 1.913895196850168
$data[15,3] = 2016
$data[15,4] = 6.97490799213798
$data[16,0] = 42321
$data[16,1] = 2015
$data[16,2] = 0.03947433952959933
$data[16,3] = 2016
$data[16,4] = -1.259568900987018
$data[17,0] = 42503
$data[17,1] = 2016
$data[17,2] = 4.861901970953975
$data[17,3] = 2017
$data[17,4] = 9.631040506010535
$data[18,0] = 42689
$data[18,1] = 2016
$data[18,2] = 2.192778679161944
$data[18,3] = 2017
$data[18,4] = 0.5033587260849126
$data[19,0] = 42867
$data[19,1] = 2017
$data[19,2] = 4.115488239647713
$data[19,3] = 2018
$data[19,4] = 9.523050046161053
$data[20,0] = 43053
$data[20,1] = 2017
$data[20,2] = 3.40836448860673
$data[20,3] = 2018
$data[20,4] = -0.3010260522302355
$data[21,0] = 43145
$data[21,1] = 2018
$data[21,2] = -0.117003051846476
$data[21,3] = 2019
$data[21,4] = -1.421675245489551
$data[22,0] = 43235
$data[22,1] = 2018
$data[22,2] = 4.07381142256642
$data[22,3] = 2019
$data[22,4] = 8.54956688663686
$data[23,0] = 43326
$data[23,1] = 2018
$data[23,2] = 2.573390224036864
$data[23,3] = 2019
$data[23,4] = 2.706783531850476
$data[24,0] = 43418
$data[24,1] = 2018
$data[24,2] = 2.799070570134488
$data[24,3] = 2019
$data[24,4] = 3.825329033908775
$data[25,0] = 43510
$data[25,1] = 2019
$data[25,2] = 4.46295719845704
$data[25,3] = 2020
$data[25,4] = 5.237111551136597
$data[26,0] = 43600
$data[26,1] = 2019
$data[26,2] = 5.264109583376908
$data[26,3] = 2020
$data[26,4] = 7.819356632099961
$data[27,0] = 43691
$data[27,1] = 2019
$data[27,2] = 3.656441317225112
$data[27,3] = 2020
$data[27,4] = -3.202316982060605
$data[28,0] = 43783
$data[28,1] = 2019
$data[28,2] = 4.195393191694419
$data[28,3] = 2020
$data[28,4] = 3.942709467505678
$data[29,0] = 43875
$data[29,1] = 2020
$data[29,2] = 1.969952049420165
$data[29,3] = 2021
$data[29,4] = 2.423519345863356
$data[30,0] = 43966
$data[30,1] = 2020
$data[30,2] = 6.942957493752444
$data[30,3] = 2021
$data[30,4] = 17.43645097609996
$data[31,0] = 44068
$data[31,1] = 2020
$data[31,2] = 2.133862376612439
$data[31,3] = 2021
$data[31,4] = -2.182504726469814
$data[32,0] = 44159
$data[32,1] = 2020
$data[32,2] = 1.666553973046048
$data[32,3] = 2021
$data[32,4] = -0.2638638106667313
$data[33,0] = 44251
$data[33,1] = 2021
$data[33,2] = -5.603015914113896
$data[33,3] = 2022
$data[33,4] = -13.74537331374777
$data[34,0] = 44341
$data[34,1] = 2021
$data[34,2] = 2.932994663878907
$data[34,3] = 2022
$data[34,4] = 2.894715150804616
$data[35,0] = 44432
$data[35,1] = 2021
$data[35,2] = 1.773820722495745
$data[35,3] = 2022
$data[35,4] = 3.694490460041355
$data[36,0] = 44525
$data[36,1] = 2021
$data[36,2] = 1.879266440112803
$data[36,3] = 2022
$data[36,4] = -0.04532879466145889
$data[37,0] = 44617
$data[37,1] = 2022
$data[37,2] = 1.312381597381518
$data[37,3] = 2023
$data[37,4] = 8.701423588527524
$data[38,0] = 44706
$data[38,1] = 2022
$data[38,2] = -0.7538332529782865
$data[38,3] = 2023
$data[38,4] = -4.784481399264983
$data[39,0] = 44798
$data[39,1] = 2022
$data[39,2] = -2.404913754290983
$data[39,3] = 2023
$data[39,4] = -1.736119732506514
$data[40,0] = 44890
$data[40,1] = 2022
$data[40,2] = -2.620683231370946
$data[40,3] = 2023
$data[40,4] = -3.179374983142691
$data[41,0] = 44981
$data[41,1] = 2023
$data[41,2] = -3.324604708321111
$data[41,3] = 2024
$data[41,4] = 2.584340770833982
$data[42,0] = 45071
$data[42,1] = 2023
$data[42,2] = -2.598185084325777
$data[42,3] = 2024
$data[42,4] = -3.071148328823314
$data[43,0] = 45163
$data[43,1] = 2023
$data[43,2] = -2.901570548279864
$data[43,3] = 2024
$data[43,4] = -3.268007849027199
$data[44,0] = 45254
$data[44,1] = 2023
$data[44,2] = -3.036556262700274
$data[44,3] = 2024
$data[44,4] = -2.51939929628594
$data[45,0] = 45345
$data[45,1] = 2024
$data[45,2] = -2.440053088416461
$data[45,3] = 2025
$data[45,4] = -0.6673203033532138
$data[46,0] = 45436
$data[46,1] = 2024
$data[46,2] = -1.287508943286542
$data[46,3] = 2025
$data[46,4] = -4.200823682253607
$data[47,0] = 45534
$data[47,1] = 2024
$data[47,2] = -2.859191689251428
$data[47,3] = 2025
$data[47,4] = -2.883845406532493
$data[48,0] = 45618
$data[48,1] = 2024
$data[48,2] = -2.953443685011514
$data[48,3] = 2025
$data[48,4] = -2.112604539331953
$data[49,0] = 45713
$data[49,1] = 2025
$data[49,2] = -0.6052121327035698
$data[49,3] = 2026
$data[49,4] = -0.8972538974235111
$data[50,0] = 45800
$data[50,1] = 2025
$data[50,2] = 0.1893861904177951
$data[50,3] = 2026
$data[50,4] = 0.1126281723122791
$ws.Range("A2:E52").Value = $data

# Append new row 53, copying the date-style format from row 52 column A
$ws.Range("A52").Copy($ws.Range("A53"))
$ws.Range("A53").Value = 45891
$ws.Range("B53").Value = 2025
$ws.Range("C53").Value = -1.131442475565558
$ws.Range("D53").Value = 2026
$ws.Range("E53").Value = -1.463023257418061

Write-Host "done"
